# Weekly refresh: a new daily price record is inserted at the top of the
# Jengibre series (row 67), pushing the existing history down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 67

# Shift rows 67:104 down to 68:105, inserting a blank row 67 in their place.
$ws.Rows.Item($newRow).Insert()

# Populate the newly inserted row with the latest market record.
$ws.Cells.Item($newRow, 1).Value = 6
$ws.Cells.Item($newRow, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($newRow, 3).Value = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value = 44830
$ws.Cells.Item($newRow, 5).Value = 13
$ws.Cells.Item($newRow, 6).Value = 100114007
$ws.Cells.Item($newRow, 7).Value = "Jengibre"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 320
$ws.Cells.Item($newRow, 11).Value = 13000
$ws.Cells.Item($newRow, 12).Value = 14000
$ws.Cells.Item($newRow, 13).Value = 13531
$ws.Cells.Item($newRow, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item($newRow, 15).Value = "Perú"
$ws.Cells.Item($newRow, 16).Value = 1041
$ws.Cells.Item($newRow, 17).Value = 13
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
